# Schedule.xlsx edit script
# Reproduces the changes described in the commit:
#  - Make "Schedule" the active/selected sheet (was "Plan")
#  - Extend the 30-minute time ladder in column B from row 30 down to row 36
#  - Fill in the Saturday (column H) activity cells for rows 29-36
#  - Update the selection on the Schedule sheet to B29:B36

$wb = $excel.ActiveWorkbook
$schedule = $wb.Worksheets.Item("Schedule")
$plan = $wb.Worksheets.Item("Plan")

# --- Column B: extend the half-hour time series down through row 36 ---
# Row 30 already holds "=B29+TIME(0, 30, 0)" (part of the shared formula
# that starts at B6). Continue the same 30-minute pattern for B31:B36,
# replacing the old 1-minute-step formulas/literal that used to live there.
$schedule.Range("B31").Formula = "=B30+TIME(0, 30, 0)"
$schedule.Range("B32").Formula = "=B31+TIME(0, 30, 0)"
$schedule.Range("B33").Formula = "=B32+TIME(0, 30, 0)"
$schedule.Range("B34").Formula = "=B33+TIME(0, 30, 0)"
$schedule.Range("B35").Formula = "=B34+TIME(0, 30, 0)"
$schedule.Range("B36").Formula = "=B35+TIME(0, 30, 0)"

# --- Column H (Saturday): set the activity names for rows 29-36 ---
$schedule.Range("H29").Value = "Оценка недели"
$schedule.Range("H30").Value = "Оценка недели"
$schedule.Range("H31").Value = "Оценка недели"
$schedule.Range("H32").Value = "Оценка недели"
$schedule.Range("H33").Value = "Биомеханика"
$schedule.Range("H34").Value = "Анализ сделок"
$schedule.Range("H35").Value = "Анализ сделок"
$schedule.Range("H36").Value = "Ужин"

# --- Make "Schedule" the active sheet/tab and set its selection ---
$schedule.Activate()
$schedule.Range("B29:B36").Select()

Write-Host "Edit complete"
